# Generate Report for Handback
# Marks the zh-cn and de-de handback rows as failed, and records the
# handback-transform error detail describing the filename mismatch.
# The same status text ("Ready for handoff" -> "Handback transform failed")
# is shared by the per-language sheets' Status column and by the Overview
# sheet's summary columns for each language, so all of those cells are
# updated together.

$wb = $excel.ActiveWorkbook

$statusText = "Handback transform failed"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusText
$zhcn.Range("K3").Value = "Handback file name: 0hs04bbe.3fn is different with handoff file name: 336e0441-8455-48ab-a0bb-de862a5e49ff.6e043e6cfbd415f0a6b91d16fdb848b052942704.zh-cn."

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusText
$dede.Range("K3").Value = "Handback file name: 0hs04bbe.3fn is different with handoff file name: 336e0441-8455-48ab-a0bb-de862a5e49ff.6e043e6cfbd415f0a6b91d16fdb848b052942704.de-de."
